$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Authors (E2) updated to the corrected full-author list; Misc. Data (J2) gets publisher
$ws.Range("E2").Value = '[Christiane Maria Prado%Jeronimo%NULL%0, Maria Eduarda Leão%Farias%NULL%0, Fernando Fonseca Almeida%Val%NULL%0, Vanderson Souza%Sampaio%NULL%0, Marcia Almeida Araújo%Alexandre%NULL%0, Gisely Cardoso%Melo%NULL%0, Izabella Picinin%Safe%NULL%0, Mayla Gabriela Silva%Borba%NULL%0, Rebeca Linhares%Abreu-Netto%NULL%0, Alex Bezerra Silva%Maciel%NULL%0, João Ricardo Silva%Neto%NULL%0, Lucas Barbosa%Oliveira%NULL%0, Erick Frota Gomes%Figueiredo%NULL%0, Kelry Mazurega Oliveira%Dinelly%NULL%0, Maria Gabriela de Almeida%Rodrigues%NULL%0, Marcelo%Brito%NULL%0, Maria Paula Gomes%Mourão%NULL%0, Guilherme Augusto%Pivoto João%NULL%0, Ludhmila Abrahão%Hajjar%NULL%0, Quique%Bassat%NULL%0, Gustavo Adolfo Sierra%Romero%NULL%0, Felipe Gomes%Naveca%NULL%0, Heline Lira%Vasconcelos%NULL%0, Michel de Araújo%Tavares%NULL%0, José Diego%Brito-Sousa%NULL%0, Fabio Trindade Maranhão%Costa%NULL%0, Maurício Lacerda%Nogueira%NULL%0, Djane%Baía-da-Silva%NULL%0, Mariana Simão%Xavier%NULL%0, Wuelton Marcelo%Monteiro%NULL%0, Marcus Vinícius Guimarães%Lacerda%marcuslacerda.br@gmail.com%0, NULL%NULL%NULL%0]'
$ws.Range("J2").Value = 'Oxford University Press'

# Row 3: Misc. Data (J3) gets publisher
$ws.Range("J3").Value = 'Massachusetts Medical Society'

# Row 4: full-text lookup failed -> reset Title/Authors/ID/ID Format to "unknown" placeholders
$ws.Range("C4").Value = 'Unknown Title'
$ws.Range("E4").Value = '[]'
$ws.Range("F4").Value = 'not found'
$ws.Range("G4").Value = 'N/A'

# Row 5: full-text lookup failed -> reset Title/Abstract/ID/ID Format to "unknown" placeholders
$ws.Range("C5").Value = 'Unknown Title'
$ws.Range("D5").Value = 'Unknown Abstract'
$ws.Range("F5").Value = 'not found'
$ws.Range("G5").Value = 'N/A'

# Row 6: full-text lookup failed -> reset Title/Abstract/Authors/ID/ID Format to "unknown" placeholders
$ws.Range("C6").Value = 'Unknown Title'
$ws.Range("D6").Value = 'Unknown Abstract'
$ws.Range("E6").Value = '[]'
$ws.Range("F6").Value = 'not found'
$ws.Range("G6").Value = 'N/A'
